$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 453
$ws.Cells.Item(453, 1).Value = 452.0
$ws.Cells.Item(453, 2).Value = 'Saturday, Jan 14'
$ws.Cells.Item(453, 3).Value = '5:40 AM'
$ws.Cells.Item(453, 4).Value = 'LO3910'
$ws.Cells.Item(453, 5).Value = 'Warsaw'
$ws.Cells.Item(453, 6).Value = '(WAW)'
$ws.Cells.Item(453, 7).Value = 'LOT '
$ws.Cells.Item(453, 8).Value = 'E190'
$ws.Cells.Item(453, 9).Value = '(SP-LMH)'
$ws.Cells.Item(453, 10).Value = '5:43 AM'
$ws.Cells.Item(453, 12).Value = '0 hours, 3 minutes'

# Row 454
$ws.Cells.Item(454, 1).Value = 453.0
$ws.Cells.Item(454, 2).Value = 'Saturday, Jan 14'
$ws.Cells.Item(454, 3).Value = '5:55 AM'
$ws.Cells.Item(454, 4).Value = 'W65091'
$ws.Cells.Item(454, 5).Value = 'Lyon'
$ws.Cells.Item(454, 6).Value = '(LYS)'
$ws.Cells.Item(454, 7).Value = 'Wizz Air '
$ws.Cells.Item(454, 8).Value = 'A321'
$ws.Cells.Item(454, 9).Value = '(HA-LXO)'
$ws.Cells.Item(454, 10).Value = '6:05 AM'
$ws.Cells.Item(454, 12).Value = '0 hours, 10 minutes'

# Row 455
$ws.Cells.Item(455, 1).Value = 454.0
$ws.Cells.Item(455, 2).Value = 'Saturday, Jan 14'
$ws.Cells.Item(455, 3).Value = '6:00 AM'
$ws.Cells.Item(455, 4).Value = 'KL1992'
$ws.Cells.Item(455, 5).Value = 'Amsterdam'
$ws.Cells.Item(455, 6).Value = '(AMS)'
$ws.Cells.Item(455, 7).Value = 'KLM '
$ws.Cells.Item(455, 8).Value = 'E295'
$ws.Cells.Item(455, 9).Value = '(PH-NXB)'
$ws.Cells.Item(455, 10).Value = '6:08 AM'
$ws.Cells.Item(455, 12).Value = '0 hours, 8 minutes'

# Row 456
$ws.Cells.Item(456, 1).Value = 455.0
$ws.Cells.Item(456, 2).Value = 'Saturday, Jan 14'
$ws.Cells.Item(456, 3).Value = '6:10 AM'
$ws.Cells.Item(456, 4).Value = 'FR6208'
$ws.Cells.Item(456, 5).Value = 'Madrid'
$ws.Cells.Item(456, 6).Value = '(MAD)'
$ws.Cells.Item(456, 7).Value = 'Buzz '
$ws.Cells.Item(456, 8).Value = 'B38M'
$ws.Cells.Item(456, 9).Value = '(SP-RZH)'
$ws.Cells.Item(456, 10).Value = '6:34 AM'
$ws.Cells.Item(456, 12).Value = '0 hours, 24 minutes'

# Row 457
$ws.Cells.Item(457, 1).Value = 456.0
$ws.Cells.Item(457, 2).Value = 'Saturday, Jan 14'
$ws.Cells.Item(457, 3).Value = '6:20 AM'
$ws.Cells.Item(457, 4).Value = 'LH1371'
$ws.Cells.Item(457, 5).Value = 'Frankfurt'
$ws.Cells.Item(457, 6).Value = '(FRA)'
$ws.Cells.Item(457, 7).Value = 'Lufthansa '
$ws.Cells.Item(457, 8).Value = 'A21N'
$ws.Cells.Item(457, 9).Value = '(D-AIEK)'
$ws.Cells.Item(457, 10).Value = '6:24 AM'
$ws.Cells.Item(457, 12).Value = '0 hours, 4 minutes'

# Row 458
$ws.Cells.Item(458, 1).Value = 457.0
$ws.Cells.Item(458, 2).Value = 'Saturday, Jan 14'
$ws.Cells.Item(458, 3).Value = '6:25 AM'
$ws.Cells.Item(458, 4).Value = 'FR6310'
$ws.Cells.Item(458, 5).Value = 'Oslo'
$ws.Cells.Item(458, 6).Value = '(TRF)'
$ws.Cells.Item(458, 7).Value = 'Ryanair '
$ws.Cells.Item(458, 8).Value = 'B738'
$ws.Cells.Item(458, 9).Value = '(SP-RKC)'
$ws.Cells.Item(458, 10).Value = '6:18 AM'
$ws.Cells.Item(458, 12).Value = '0 hours, -7 minutes'

# Row 459
$ws.Cells.Item(459, 1).Value = 458.0
$ws.Cells.Item(459, 2).Value = 'Saturday, Jan 14'
$ws.Cells.Item(459, 3).Value = '6:25 AM'
$ws.Cells.Item(459, 4).Value = 'W65097'
$ws.Cells.Item(459, 5).Value = 'Tel Aviv'
$ws.Cells.Item(459, 6).Value = '(TLV)'
$ws.Cells.Item(459, 7).Value = 'Wizz Air '
$ws.Cells.Item(459, 8).Value = 'A21N'
$ws.Cells.Item(459, 9).Value = '(HA-LVO)'
$ws.Cells.Item(459, 10).Value = '6:40 AM'
$ws.Cells.Item(459, 12).Value = '0 hours, 15 minutes'

# Row 460
$ws.Cells.Item(460, 1).Value = 459.0
$ws.Cells.Item(460, 2).Value = 'Saturday, Jan 14'
$ws.Cells.Item(460, 3).Value = '6:50 AM'
$ws.Cells.Item(460, 4).Value = 'FR6228'
$ws.Cells.Item(460, 5).Value = 'Tel Aviv'
$ws.Cells.Item(460, 6).Value = '(TLV)'
$ws.Cells.Item(460, 7).Value = 'Ryanair '
$ws.Cells.Item(460, 8).Value = 'B38M'
$ws.Cells.Item(460, 9).Value = '(SP-RZK)'
$ws.Cells.Item(460, 10).Value = '6:58 AM'
$ws.Cells.Item(460, 12).Value = '0 hours, 8 minutes'

# Row 461
$ws.Cells.Item(461, 1).Value = 460.0
$ws.Cells.Item(461, 2).Value = 'Saturday, Jan 14'
$ws.Cells.Item(461, 3).Value = '6:50 AM'
$ws.Cells.Item(461, 4).Value = 'FR6322'
$ws.Cells.Item(461, 5).Value = 'Amman'
$ws.Cells.Item(461, 6).Value = '(AMM)'
$ws.Cells.Item(461, 7).Value = 'Buzz '
$ws.Cells.Item(461, 8).Value = 'B38M'
$ws.Cells.Item(461, 9).Value = '(SP-RZB)'
$ws.Cells.Item(461, 10).Value = '7:01 AM'
$ws.Cells.Item(461, 12).Value = '0 hours, 11 minutes'

# Row 462
$ws.Cells.Item(462, 1).Value = 461.0
$ws.Cells.Item(462, 2).Value = 'Saturday, Jan 14'
$ws.Cells.Item(462, 3).Value = '7:10 AM'
$ws.Cells.Item(462, 4).Value = 'FR3680'
$ws.Cells.Item(462, 5).Value = 'Birmingham'
$ws.Cells.Item(462, 6).Value = '(BHX)'
$ws.Cells.Item(462, 7).Value = 'Buzz '
$ws.Cells.Item(462, 8).Value = 'B38M'
$ws.Cells.Item(462, 9).Value = '(SP-RZD)'
$ws.Cells.Item(462, 10).Value = '7:18 AM'
$ws.Cells.Item(462, 12).Value = '0 hours, 8 minutes'

# Row 463
$ws.Cells.Item(463, 1).Value = 462.0
$ws.Cells.Item(463, 2).Value = 'Saturday, Jan 14'
$ws.Cells.Item(463, 3).Value = '7:15 AM'
$ws.Cells.Item(463, 4).Value = 'FR5906'
$ws.Cells.Item(463, 5).Value = 'Turin'
$ws.Cells.Item(463, 6).Value = '(TRN)'
$ws.Cells.Item(463, 7).Value = 'Buzz '
$ws.Cells.Item(463, 8).Value = 'B38M'
$ws.Cells.Item(463, 9).Value = '(SP-RZF)'
$ws.Cells.Item(463, 10).Value = '7:23 AM'
$ws.Cells.Item(463, 12).Value = '0 hours, 8 minutes'

# Row 464
$ws.Cells.Item(464, 1).Value = 463.0
$ws.Cells.Item(464, 2).Value = 'Saturday, Jan 14'
$ws.Cells.Item(464, 3).Value = '7:25 AM'
$ws.Cells.Item(464, 4).Value = 'OS600'
$ws.Cells.Item(464, 5).Value = 'Vienna'
$ws.Cells.Item(464, 6).Value = '(VIE)'
$ws.Cells.Item(464, 7).Value = 'Austrian Airlines '
$ws.Cells.Item(464, 8).Value = 'E195'
$ws.Cells.Item(464, 9).Value = '(OE-LWP)'
$ws.Cells.Item(464, 10).Value = '7:21 AM'
$ws.Cells.Item(464, 12).Value = '0 hours, -4 minutes'

# Row 465
$ws.Cells.Item(465, 1).Value = 464.0
$ws.Cells.Item(465, 2).Value = 'Saturday, Jan 14'
$ws.Cells.Item(465, 3).Value = '7:25 AM'
$ws.Cells.Item(465, 4).Value = 'W65043'
$ws.Cells.Item(465, 5).Value = 'Stavanger'
$ws.Cells.Item(465, 6).Value = '(SVG)'
$ws.Cells.Item(465, 7).Value = 'Wizz Air '
$ws.Cells.Item(465, 8).Value = 'A21N'
$ws.Cells.Item(465, 9).Value = '(HA-LZI)'
$ws.Cells.Item(465, 10).Value = '7:31 AM'
$ws.Cells.Item(465, 12).Value = '0 hours, 6 minutes'

# Row 466
$ws.Cells.Item(466, 1).Value = 465.0
$ws.Cells.Item(466, 2).Value = 'Saturday, Jan 14'
$ws.Cells.Item(466, 3).Value = '7:30 AM'
$ws.Cells.Item(466, 4).Value = 'LH1625'
$ws.Cells.Item(466, 5).Value = 'Munich'
$ws.Cells.Item(466, 6).Value = '(MUC)'
$ws.Cells.Item(466, 7).Value = 'Lufthansa '
$ws.Cells.Item(466, 8).Value = 'A320'
$ws.Cells.Item(466, 9).Value = '(D-AIWA)'
$ws.Cells.Item(466, 10).Value = '7:34 AM'
$ws.Cells.Item(466, 12).Value = '0 hours, 4 minutes'

# Row 467
$ws.Cells.Item(467, 1).Value = 466.0
$ws.Cells.Item(467, 2).Value = 'Saturday, Jan 14'
$ws.Cells.Item(467, 3).Value = '7:40 AM'
$ws.Cells.Item(467, 4).Value = 'FR3036'
$ws.Cells.Item(467, 5).Value = 'Barcelona'
$ws.Cells.Item(467, 6).Value = '(BCN)'
$ws.Cells.Item(467, 7).Value = 'Ryanair '
$ws.Cells.Item(467, 8).Value = 'B738'
$ws.Cells.Item(467, 9).Value = '(SP-RKB)'
$ws.Cells.Item(467, 10).Value = '7:48 AM'
$ws.Cells.Item(467, 12).Value = '0 hours, 8 minutes'

# Row 468
$ws.Cells.Item(468, 1).Value = 467.0
$ws.Cells.Item(468, 2).Value = 'Saturday, Jan 14'
$ws.Cells.Item(468, 3).Value = '7:45 AM'
$ws.Cells.Item(468, 4).Value = 'FR6359'
$ws.Cells.Item(468, 5).Value = 'Liverpool'
$ws.Cells.Item(468, 6).Value = '(LPL)'
$ws.Cells.Item(468, 7).Value = 'Buzz '
$ws.Cells.Item(468, 8).Value = 'B38M'
$ws.Cells.Item(468, 9).Value = '(SP-RZA)'
$ws.Cells.Item(468, 10).Value = '7:50 AM'
$ws.Cells.Item(468, 12).Value = '0 hours, 5 minutes'

# Row 469
$ws.Cells.Item(469, 1).Value = 468.0
$ws.Cells.Item(469, 2).Value = 'Saturday, Jan 14'
$ws.Cells.Item(469, 3).Value = '8:20 AM'
$ws.Cells.Item(469, 4).Value = 'W65093'
$ws.Cells.Item(469, 5).Value = 'Abu Dhabi'
$ws.Cells.Item(469, 6).Value = '(AUH)'
$ws.Cells.Item(469, 7).Value = 'Wizz Air '
$ws.Cells.Item(469, 8).Value = 'A21N'
$ws.Cells.Item(469, 9).Value = '(HA-LVG)'
$ws.Cells.Item(469, 10).Value = '8:27 AM'
$ws.Cells.Item(469, 12).Value = '0 hours, 7 minutes'

# Row 470
$ws.Cells.Item(470, 1).Value = 469.0
$ws.Cells.Item(470, 2).Value = 'Saturday, Jan 14'
$ws.Cells.Item(470, 3).Value = '8:25 AM'
$ws.Cells.Item(470, 4).Value = 'FR9663'
$ws.Cells.Item(470, 5).Value = 'Rome'
$ws.Cells.Item(470, 6).Value = '(CIA)'
$ws.Cells.Item(470, 7).Value = 'Ryanair '
$ws.Cells.Item(470, 8).Value = 'B738'
$ws.Cells.Item(470, 9).Value = '(9H-QAS)'
$ws.Cells.Item(470, 10).Value = '8:31 AM'
$ws.Cells.Item(470, 12).Value = '0 hours, 6 minutes'

# Row 471
$ws.Cells.Item(471, 1).Value = 470.0
$ws.Cells.Item(471, 2).Value = 'Saturday, Jan 14'
$ws.Cells.Item(471, 3).Value = '8:45 AM'
$ws.Cells.Item(471, 4).Value = 'AY1162'
$ws.Cells.Item(471, 5).Value = 'Helsinki'
$ws.Cells.Item(471, 6).Value = '(HEL)'
$ws.Cells.Item(471, 7).Value = 'Finnair '
$ws.Cells.Item(471, 8).Value = 'E190'
$ws.Cells.Item(471, 9).Value = '(OH-LKH)'
$ws.Cells.Item(471, 10).Value = '8:44 AM'
$ws.Cells.Item(471, 12).Value = '0 hours, -1 minutes'

# Row 472
$ws.Cells.Item(472, 1).Value = 471.0
$ws.Cells.Item(472, 2).Value = 'Saturday, Jan 14'
$ws.Cells.Item(472, 3).Value = '8:50 AM'
$ws.Cells.Item(472, 4).Value = 'LO3904'
$ws.Cells.Item(472, 5).Value = 'Warsaw'
$ws.Cells.Item(472, 6).Value = '(WAW)'
$ws.Cells.Item(472, 7).Value = 'LOT (Star Alliance Livery) '
$ws.Cells.Item(472, 8).Value = 'E75S'
$ws.Cells.Item(472, 9).Value = '(SP-LIO)'
$ws.Cells.Item(472, 10).Value = '8:53 AM'
$ws.Cells.Item(472, 12).Value = '0 hours, 3 minutes'

# Row 473
$ws.Cells.Item(473, 1).Value = 472.0
$ws.Cells.Item(473, 2).Value = 'Saturday, Jan 14'
$ws.Cells.Item(473, 3).Value = '9:05 AM'
$ws.Cells.Item(473, 4).Value = 'FR6314'
$ws.Cells.Item(473, 5).Value = 'Marseille'
$ws.Cells.Item(473, 6).Value = '(MRS)'
$ws.Cells.Item(473, 7).Value = 'Ryanair '
$ws.Cells.Item(473, 8).Value = 'B738'
$ws.Cells.Item(473, 9).Value = '(9H-QAC)'
$ws.Cells.Item(473, 10).Value = '9:11 AM'
$ws.Cells.Item(473, 12).Value = '0 hours, 6 minutes'

# Row 474
$ws.Cells.Item(474, 1).Value = 473.0
$ws.Cells.Item(474, 2).Value = 'Saturday, Jan 14'
$ws.Cells.Item(474, 3).Value = '9:45 AM'
$ws.Cells.Item(474, 4).Value = 'E47903'
$ws.Cells.Item(474, 5).Value = 'Antalya'
$ws.Cells.Item(474, 6).Value = '(AYT)'
$ws.Cells.Item(474, 7).Value = 'Enter Air '
$ws.Cells.Item(474, 8).Value = 'B738'
$ws.Cells.Item(474, 9).Value = '(SP-ESH)'
$ws.Cells.Item(474, 10).Value = '9:50 AM'
$ws.Cells.Item(474, 12).Value = '0 hours, 5 minutes'

# Row 475
$ws.Cells.Item(475, 1).Value = 474.0
$ws.Cells.Item(475, 2).Value = 'Saturday, Jan 14'
$ws.Cells.Item(475, 3).Value = '10:15 AM'
$ws.Cells.Item(475, 4).Value = 'RR9501'
$ws.Cells.Item(475, 5).Value = 'Tel Aviv'
$ws.Cells.Item(475, 6).Value = '(TLV)'
$ws.Cells.Item(475, 7).Value = 'Ryanair '
$ws.Cells.Item(475, 8).Value = 'B738'
$ws.Cells.Item(475, 9).Value = '(SP-RSH)'
$ws.Cells.Item(475, 10).Value = '10:18 AM'
$ws.Cells.Item(475, 12).Value = '0 hours, 3 minutes'
